$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5579.8
$ws.Range("I19").Value = 5724.75
$ws.Range("K19").Value = 5724.75
$ws.Range("M19").Value = -5549.75

$ws.Range("H40").Value = 3227
$ws.Range("J40").Value = 3033.9167
$ws.Range("L40").Value = 3033.9167
$ws.Range("N40").Value = -3383.9167

$ws.Range("H106").Value = 743.2222
$ws.Range("I106").Value = 743.2222
$ws.Range("K106").Value = 743.2222
$ws.Range("M106").Value = -112.2222

$ws.Range("H107").Value = 558.9167
$ws.Range("I107").Value = 565.4545000000001
$ws.Range("K107").Value = 565.4545000000001
$ws.Range("M107").Value = 1354.5455

$ws.Range("H113").Value = 3561.6667
$ws.Range("I113").Value = 3799
$ws.Range("K113").Value = 3799
$ws.Range("M113").Value = -545

$ws.Range("H116").Value = 97566.86
$ws.Range("I116").Value = 187656.33
$ws.Range("J116").Value = 29999.75
$ws.Range("K116").Value = 187656.33
$ws.Range("L116").Value = 29999.75
$ws.Range("M116").Value = -184214.33
$ws.Range("N116").Value = -36883.75

$ws.Range("H132").Value = 2538.353
$ws.Range("I132").Value = 1059.325
$ws.Range("K132").Value = 3177.975
$ws.Range("M132").Value = -647.9750000000004

$ws.Range("H137").Value = 62503064
$ws.Range("J137").Value = 3463.8462
$ws.Range("L137").Value = 10391.5386
$ws.Range("N137").Value = -15491.5386

$ws.Range("H138").Value = 3123.2954
$ws.Range("I138").Value = 2037.2307
$ws.Range("J138").Value = 3578.742
$ws.Range("K138").Value = 6111.6921
$ws.Range("L138").Value = 10736.226
$ws.Range("M138").Value = -971.6921000000002
$ws.Range("N138").Value = -21016.226

$ws.Range("H141").Value = 944
$ws.Range("I141").Value = 944
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2832
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2348
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1299.6459
$ws.Range("I32").Value = 1299.6459
$ws.Range("K32").Value = 1299.6459
$ws.Range("M32").Value = -1012.6459

$ws.Range("H45").Value = 1797.3462
$ws.Range("I45").Value = 1740.3334
$ws.Range("K45").Value = 1740.3334
$ws.Range("M45").Value = -1363.3334

$ws.Range("H61").Value = 1443.1333
$ws.Range("I61").Value = 1250.878
$ws.Range("J61").Value = 3413.75
$ws.Range("K61").Value = 1250.878
$ws.Range("L61").Value = 3413.75
$ws.Range("M61").Value = -1038.878
$ws.Range("N61").Value = -3837.75

$ws.Range("H74").Value = 1962.75
$ws.Range("I74").Value = 1323
$ws.Range("J74").Value = 9000
$ws.Range("K74").Value = 1323
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -449
$ws.Range("N74").Value = -10748

$ws.Range("H77").Value = 1962.75
$ws.Range("I77").Value = 1323
$ws.Range("J77").Value = 9000
$ws.Range("K77").Value = 6615
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -2247
$ws.Range("N77").Value = -53736

$ws.Range("H95").Value = 59966.4
$ws.Range("J95").Value = 59966.4
$ws.Range("L95").Value = 59966.4
$ws.Range("N95").Value = -65458.4

$ws.Range("H110").Value = 1779.1538
$ws.Range("I110").Value = 1454.1428
$ws.Range("J110").Value = 2158.3333
$ws.Range("K110").Value = 1454.1428
$ws.Range("L110").Value = 2158.3333
$ws.Range("M110").Value = 590.8571999999999
$ws.Range("N110").Value = -6248.3333

$ws.Range("H122").Value = 1947.5
$ws.Range("I122").Value = 1430.3334
$ws.Range("K122").Value = 4291.0002
$ws.Range("M122").Value = -1841.0002

$ws.Range("H132").Value = 4453
$ws.Range("I132").Value = 4506.353
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 13519.059
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -10989.059
$ws.Range("N132").Value = -17058.5

$ws.Range("H136").Value = 1443.1333
$ws.Range("I136").Value = 1250.878
$ws.Range("J136").Value = 3413.75
$ws.Range("K136").Value = 3752.634
$ws.Range("L136").Value = 10241.25
$ws.Range("M136").Value = -1202.634
$ws.Range("N136").Value = -15341.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2093.3076
$ws.Range("I134").Value = 1350.2245
$ws.Range("K134").Value = 4050.6735
$ws.Range("M134").Value = -1515.6735

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 23333.334
$ws.Range("I51").Value = 23333.334
$ws.Range("K51").Value = 23333.334
$ws.Range("M51").Value = -22597.334

$ws.Range("H53").Value = 39999
$ws.Range("J53").Value = 39999
$ws.Range("L53").Value = 39999
$ws.Range("N53").Value = -41213

$ws.Range("H61").Value = 23333.334
$ws.Range("I61").Value = 23333.334
$ws.Range("K61").Value = 23333.334
$ws.Range("M61").Value = -22985.334

$ws.Range("H93").Value = 49999.25
$ws.Range("I93").Value = 49998.5
$ws.Range("K93").Value = 49998.5
$ws.Range("M93").Value = -48126.5

$ws.Range("H105").Value = 1214.1613
$ws.Range("I105").Value = 962.9474
$ws.Range("J105").Value = 1611.9166
$ws.Range("K105").Value = 962.9474
$ws.Range("L105").Value = 1611.9166
$ws.Range("M105").Value = 784.0526
$ws.Range("N105").Value = -5105.9166

$ws.Range("H132").Value = 133334664
$ws.Range("I132").Value = 142858480
$ws.Range("K132").Value = 428575440
$ws.Range("M132").Value = -428572910

$ws.Range("H140").Value = 92864.27
$ws.Range("J140").Value = 92864.27
$ws.Range("L140").Value = 92864.27
$ws.Range("N140").Value = -103224.27

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 121.1875
$ws.Range("I23").Value = 99.55556
$ws.Range("J23").Value = 149
$ws.Range("K23").Value = 298.66668
$ws.Range("L23").Value = 447
$ws.Range("M23").Value = -63.66667999999999
$ws.Range("N23").Value = -917

$ws.Range("H86").Value = 325.5
$ws.Range("I86").Value = 366.33334
$ws.Range("J86").Value = 203
$ws.Range("K86").Value = 1099.00002
$ws.Range("L86").Value = 609
$ws.Range("M86").Value = 86.99998000000005
$ws.Range("N86").Value = -2981

$ws.Range("H89").Value = 325.5
$ws.Range("I89").Value = 366.33334
$ws.Range("J89").Value = 203
$ws.Range("K89").Value = 3297.00006
$ws.Range("L89").Value = 1827
$ws.Range("M89").Value = 2630.99994
$ws.Range("N89").Value = -13683

$ws.Range("H113").Value = 4056.8572
$ws.Range("J113").Value = 4399.8335
$ws.Range("L113").Value = 13199.5005
$ws.Range("N113").Value = -17539.5005

$ws.Range("H139").Value = 73842
$ws.Range("I139").Value = 79214.46000000001
$ws.Range("J139").Value = 4000
$ws.Range("K139").Value = 237643.38
$ws.Range("L139").Value = 12000
$ws.Range("M139").Value = -232503.38
$ws.Range("N139").Value = -22280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 52659388
$ws.Range("I132").Value = 66686092
$ws.Range("J132").Value = 59252.75
$ws.Range("K132").Value = 200058276
$ws.Range("L132").Value = 177758.25
$ws.Range("M132").Value = -200055746
$ws.Range("N132").Value = -182818.25

$ws.Range("H136").Value = 21137.5
$ws.Range("J136").Value = 21137.5
$ws.Range("L136").Value = 63412.5
$ws.Range("N136").Value = -68512.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7035.875
$ws.Range("I16").Value = 10608.2
$ws.Range("K16").Value = 10608.2
$ws.Range("M16").Value = -10438.2

$ws.Range("H22").Value = 661.6667
$ws.Range("J22").Value = 724.625
$ws.Range("L22").Value = 724.625
$ws.Range("N22").Value = -1314.625

$ws.Range("H27").Value = 661.6667
$ws.Range("J27").Value = 724.625
$ws.Range("L27").Value = 724.625
$ws.Range("N27").Value = -938.625

$ws.Range("H40").Value = 2977.5715
$ws.Range("I40").Value = 2977.5715
$ws.Range("K40").Value = 2977.5715
$ws.Range("M40").Value = -2841.5715

$ws.Range("H132").Value = 2848.6875
$ws.Range("I132").Value = 2838.2666
$ws.Range("K132").Value = 8514.799800000001
$ws.Range("M132").Value = -5984.799800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 69473
$ws.Range("J99").Value = 69473
$ws.Range("L99").Value = 69473
$ws.Range("N99").Value = -75463

$ws.Range("H122").Value = 1408.9412
$ws.Range("I122").Value = 1408.9412
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4226.8236
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1776.8236
$ws.Range("N122").Value = ""

$ws.Range("H126").Value = 2229.4482
$ws.Range("I126").Value = 1689.5883
$ws.Range("K126").Value = 5068.7649
$ws.Range("M126").Value = -2598.7649

$ws.Range("H132").Value = 50009524
$ws.Range("I132").Value = 100003050
$ws.Range("J132").Value = 16000
$ws.Range("K132").Value = 300009150
$ws.Range("L132").Value = 48000
$ws.Range("M132").Value = -300006620
$ws.Range("N132").Value = -53060

$ws.Range("I136").Value = 4977026
$ws.Range("K136").Value = 14931078
$ws.Range("M136").Value = -14928528
